$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MaDon (col A) and SoTien (col E) are cleared for the first two data rows ...
$ws.Range("A2:A3").ClearContents()
$ws.Range("E2:E3").ClearContents()

# ... and all remaining data rows (4-7) are cleared completely.
$ws.Range("A4:F7").ClearContents()
